# Fill in the two "System Requirement" cells for the "Volunteer Hours Log
# Page with Mobile Capability" row (row 10) on Sheet1: these describe the
# account-creation and time/date/location logging requirements.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# F10 is set first so it lands on shared-string index 34 and E10 on 35,
# matching the authored workbook.
$ws.Range("F10").Value = "Users will use their accounts to log specific time, date, locations, and tasks performed.  A Smithgall Woods representative will be able to log in and confirm tasks/hours logged.  "
$ws.Range("E10").Value = "Volunteers will create accounts in which name, date of birth, home address, email, phone numbers will be recorded.  Users will also create a password for their account.  Once accounts are created, a Smithgall Woods representative will validate user by logging into system.  "
